# Weekly update: insert a new week's worth of Mango price data
# (Especial / Primera / Segunda) at the top of the dated price table,
# shifting the existing rows 547:578 down to 550:581.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at row 547 (formats are inherited from the row above,
# matching how the rest of the table is formatted - only column D carries
# the date style).
$ws.Rows("547:549").Insert()

$grades = "Especial", "Primera", "Segunda"

for ($i = 0; $i -lt 3; $i++) {
    $r = 547 + $i

    $ws.Cells.Item($r, 1).Value = 8
    $ws.Cells.Item($r, 2).Value = "Terminal La Palmera de La Serena"
    $ws.Cells.Item($r, 3).Value = "Coquimbo"
    $ws.Cells.Item($r, 4).Value = 44585
    $ws.Cells.Item($r, 5).Value = 4
    $ws.Cells.Item($r, 6).Value = "Fruta"
    $ws.Cells.Item($r, 7).Value = 100108
    $ws.Cells.Item($r, 8).Value = "Tropicales y subtropicales"
    $ws.Cells.Item($r, 9).Value = 100108002
    $ws.Cells.Item($r, 10).Value = "Mango"
    $ws.Cells.Item($r, 11).Value = "Sin especificar"
    $ws.Cells.Item($r, 12).Value = $grades[$i]
    $ws.Cells.Item($r, 13).Value = 512
    $ws.Cells.Item($r, 14).Value = 7500
    $ws.Cells.Item($r, 15).Value = 8000
    $ws.Cells.Item($r, 16).Value = 7750
    $ws.Cells.Item($r, 17).Value = "`$/bandeja 4 kilos"
    $ws.Cells.Item($r, 18).Value = "Perú"
    $ws.Cells.Item($r, 19).Value = 1938
    $ws.Cells.Item($r, 20).Value = 4
}
